$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A ("Group") before the existing data, shifting
# Material..Source from A:F to B:G.
$ws.Columns.Item(1).Insert()

# Header for the new column.
$ws.Cells.Item(1, 1).Value = "Group"

# Group labels for each material row (rows 2-7). Written in A, B, C, O,
# Cop, F order so new shared-string entries are interned in that sequence.
$ws.Cells.Item(2, 1).Value = "A"
$ws.Cells.Item(3, 1).Value = "B"
$ws.Cells.Item(5, 1).Value = "C"
$ws.Cells.Item(4, 1).Value = "O"
$ws.Cells.Item(6, 1).Value = "Cop"
$ws.Cells.Item(7, 1).Value = "F"

# Style the new Group column cells (rows 2-7) with a distinct fill and a
# thin box border on all sides.
$groupRange = $ws.Range("A2:A7")
$groupRange.Borders.LineStyle = 1
$groupRange.Interior.ThemeColor = 1
$groupRange.Interior.TintAndShade = -0.34998626667073579

# Column width for the new Group column.
$ws.Columns.Item(1).ColumnWidth = 6.42578125

# Update selection to match target state.
$ws.Range("D13").Select()
